# Update the shipping/instance lookup sheet:
#   - Column A header "Region" -> "ID"
#   - Column A values change from AWS availability-zone names to EC2
#     instance ids
#   - Column A is widened to fit the new, longer values
#   - Selection moves to A3 (last cell touched interactively)
#
# Note: column B ("Shipping_Fee" header + numeric fee values) is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the instance-id values first, then the header, so the shared-string
# table is built in the same order as the authored workbook
# (Shipping_Fee, i-0a64cd90fd5757c27 , i-0db5e5e4af2bd0772,
#  i-0dd550b7fbd9ae483, ID).
$ws.Range("A2").Value = "i-0a64cd90fd5757c27 "
$ws.Range("A3").Value = "i-0db5e5e4af2bd0772"
$ws.Range("A4").Value = "i-0dd550b7fbd9ae483"
$ws.Range("A1").Value = "ID"

# Widen column A to fit the new (longer) id strings - resulting stored
# column width of 20 characters.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(1).ColumnWidth = 19.17

# Leave the final selection on A3, matching the saved cursor position.
$ws.Range("A3").Select() | Out-Null
